$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 27 de Marzo de 2020 a las 09:12'
$ws.Cells.Item(11, 5).Value = 11486
$ws.Cells.Item(11, 7).Value = 2
$ws.Cells.Item(11, 8).Value = 194
$ws.Cells.Item(25, 1).Value = 'Chequia'
$ws.Cells.Item(25, 2).Value = 2062
$ws.Cells.Item(25, 3).Value = 137
$ws.Cells.Item(25, 4).Value = 10
$ws.Cells.Item(25, 5).Value = 2043
$ws.Cells.Item(25, 6).Value = 34
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 9
$ws.Cells.Item(26, 1).Value = 'Malasia'
$ws.Cells.Item(26, 2).Value = 2031
$ws.Cells.Item(26, 4).Value = 215
$ws.Cells.Item(26, 5).Value = 1792
$ws.Cells.Item(26, 6).Value = 45
$ws.Cells.Item(26, 7).Value = 1
$ws.Cells.Item(26, 8).Value = 24
$ws.Cells.Item(53, 1).Value = 'Croacia'
$ws.Cells.Item(53, 2).Value = 551
$ws.Cells.Item(53, 3).Value = 56
$ws.Cells.Item(53, 4).Value = 37
$ws.Cells.Item(53, 5).Value = 511
$ws.Cells.Item(53, 6).Value = 14
$ws.Cells.Item(53, 8).Value = 3
$ws.Cells.Item(54, 1).Value = 'Catar'
$ws.Cells.Item(54, 2).Value = 549
$ws.Cells.Item(54, 4).Value = 43
$ws.Cells.Item(54, 5).Value = 506
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(55, 1).Value = 'Estonia'
$ws.Cells.Item(55, 2).Value = 538
$ws.Cells.Item(55, 4).Value = 8
$ws.Cells.Item(55, 5).Value = 529
$ws.Cells.Item(55, 6).Value = 6
$ws.Cells.Item(55, 8).Value = 1
$ws.Cells.Item(59, 4).Value = 210
$ws.Cells.Item(59, 5).Value = 244
$ws.Cells.Item(70, 1).Value = 'Letonia'
$ws.Cells.Item(70, 2).Value = 280
$ws.Cells.Item(70, 3).Value = 36
$ws.Cells.Item(70, 4).Value = 1
$ws.Cells.Item(70, 5).Value = 279
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(71, 1).Value = 'Bulgaria'
$ws.Cells.Item(71, 2).Value = 276
$ws.Cells.Item(71, 3).Value = 12
$ws.Cells.Item(71, 5).Value = 265
$ws.Cells.Item(71, 6).Value = 8
$ws.Cells.Item(71, 8).Value = 3
$ws.Cells.Item(72, 1).Value = 'Marruecos'
$ws.Cells.Item(72, 2).Value = 275
$ws.Cells.Item(72, 3).Value = 0
$ws.Cells.Item(72, 4).Value = 8
$ws.Cells.Item(72, 5).Value = 256
$ws.Cells.Item(72, 6).Value = 1
$ws.Cells.Item(72, 8).Value = 11
$ws.Cells.Item(73, 1).Value = 'Taiwan'
$ws.Cells.Item(73, 2).Value = 267
$ws.Cells.Item(73, 3).Value = 15
$ws.Cells.Item(73, 4).Value = 30
$ws.Cells.Item(73, 5).Value = 235
$ws.Cells.Item(73, 8).Value = 2
$ws.Cells.Item(78, 1).Value = 'Ucrania'
$ws.Cells.Item(78, 2).Value = 218
$ws.Cells.Item(78, 3).Value = 22
$ws.Cells.Item(78, 4).Value = 4
$ws.Cells.Item(78, 5).Value = 209
$ws.Cells.Item(78, 8).Value = 5
$ws.Cells.Item(79, 1).Value = 'Jordania'
$ws.Cells.Item(79, 2).Value = 212
$ws.Cells.Item(79, 4).Value = 1
$ws.Cells.Item(79, 5).Value = 211
$ws.Cells.Item(79, 6).Value = 0
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(80, 1).Value = 'San Marino'
$ws.Cells.Item(80, 4).Value = 4
$ws.Cells.Item(80, 5).Value = 183
$ws.Cells.Item(80, 6).Value = 12
$ws.Cells.Item(80, 8).Value = 21
$ws.Cells.Item(81, 1).Value = 'Kuwait'
$ws.Cells.Item(81, 2).Value = 208
$ws.Cells.Item(81, 4).Value = 49
$ws.Cells.Item(81, 5).Value = 159
$ws.Cells.Item(81, 6).Value = 7
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(82, 1).Value = 'Republica de Macedonia'
$ws.Cells.Item(82, 2).Value = 201
$ws.Cells.Item(82, 4).Value = 3
$ws.Cells.Item(82, 5).Value = 195
$ws.Cells.Item(82, 6).Value = 1
$ws.Cells.Item(82, 8).Value = 3
$ws.Cells.Item(83, 1).Value = 'Tunez'
$ws.Cells.Item(83, 2).Value = 197
$ws.Cells.Item(83, 4).Value = 2
$ws.Cells.Item(83, 6).Value = 10
$ws.Cells.Item(84, 2).Value = 192
$ws.Cells.Item(84, 3).Value = 1
$ws.Cells.Item(84, 4).Value = 5
$ws.Cells.Item(84, 5).Value = 184
$ws.Cells.Item(85, 5).Value = 173
$ws.Cells.Item(85, 7).Value = 1
$ws.Cells.Item(85, 8).Value = 2
$ws.Cells.Item(90, 4).Value = 54
$ws.Cells.Item(90, 5).Value = 86
$ws.Cells.Item(90, 6).Value = 2
$ws.Cells.Item(118, 1).Value = 'Kirguistan'
$ws.Cells.Item(118, 2).Value = 58
$ws.Cells.Item(118, 3).Value = 14
$ws.Cells.Item(118, 5).Value = 58
$ws.Cells.Item(119, 1).Value = 'Liechtenstein'
$ws.Cells.Item(119, 2).Value = 56
$ws.Cells.Item(119, 3).Value = 0
$ws.Cells.Item(119, 4).Value = 0
$ws.Cells.Item(119, 5).Value = 56
$ws.Cells.Item(119, 6).Value = 0
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(120, 1).Value = 'Paraguay'
$ws.Cells.Item(120, 2).Value = 52
$ws.Cells.Item(120, 3).Value = 11
$ws.Cells.Item(120, 4).Value = 1
$ws.Cells.Item(120, 5).Value = 48
$ws.Cells.Item(120, 6).Value = 1
$ws.Cells.Item(121, 1).Value = 'Consejo Danes para los Refugiados'
$ws.Cells.Item(121, 2).Value = 51
$ws.Cells.Item(121, 4).Value = 2
$ws.Cells.Item(121, 5).Value = 46
$ws.Cells.Item(121, 8).Value = 3
$ws.Cells.Item(122, 1).Value = 'Ruanda'
$ws.Cells.Item(122, 2).Value = 50
$ws.Cells.Item(122, 3).Value = 0
$ws.Cells.Item(122, 4).Value = 0
$ws.Cells.Item(122, 5).Value = 50
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(123, 1).Value = 'Banglades'
$ws.Cells.Item(123, 2).Value = 48
$ws.Cells.Item(123, 3).Value = 4
$ws.Cells.Item(123, 4).Value = 11
$ws.Cells.Item(123, 5).Value = 32
$ws.Cells.Item(123, 6).Value = 1
$ws.Cells.Item(123, 8).Value = 5
$ws.Cells.Item(143, 1).Value = 'Nueva Caledonia'
$ws.Cells.Item(143, 3).Value = 1
$ws.Cells.Item(143, 4).Value = 0
$ws.Cells.Item(143, 5).Value = 15
$ws.Cells.Item(144, 1).Value = 'Bermudas'
$ws.Cells.Item(144, 2).Value = 15
$ws.Cells.Item(144, 4).Value = 2
$ws.Cells.Item(144, 5).Value = 13
$ws.Cells.Item(161, 1).Value = 'Antigua y Barbuda'
$ws.Cells.Item(162, 1).Value = 'Seychelles'
$ws.Cells.Item(163, 1).Value = 'Mozambique'
$ws.Cells.Item(167, 1).Value = 'Laos'
$ws.Cells.Item(170, 1).Value = 'Eritrea'
$ws.Cells.Item(172, 1).Value = 'Siria'
$ws.Cells.Item(173, 1).Value = 'Montserrat'
$ws.Cells.Item(176, 1).Value = 'Zimbabue'
$ws.Cells.Item(176, 3).Value = 2
$ws.Cells.Item(177, 1).Value = 'Guyana'
$ws.Cells.Item(177, 2).Value = 5
$ws.Cells.Item(177, 8).Value = 1
$ws.Cells.Item(179, 1).Value = 'Angola'
$ws.Cells.Item(180, 1).Value = 'Santa Sede'
$ws.Cells.Item(181, 1).Value = 'Congo'
$ws.Cells.Item(182, 1).Value = 'Guinea'
$ws.Cells.Item(182, 2).Value = 4
$ws.Cells.Item(182, 5).Value = 4
$ws.Cells.Item(183, 1).Value = 'Liberia'
$ws.Cells.Item(184, 1).Value = 'Mauritania'
$ws.Cells.Item(185, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(186, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(187, 1).Value = 'San Bartolome'
$ws.Cells.Item(188, 1).Value = 'Republica del Chad'
$ws.Cells.Item(188, 3).Value = 0
$ws.Cells.Item(189, 1).Value = 'Butan'
$ws.Cells.Item(189, 3).Value = 1
$ws.Cells.Item(189, 4).Value = 0
$ws.Cells.Item(189, 5).Value = 3
$ws.Cells.Item(190, 1).Value = 'Santa Lucia'
$ws.Cells.Item(190, 4).Value = 1
$ws.Cells.Item(190, 8).Value = 0
$ws.Cells.Item(191, 1).Value = 'Sudan'
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 8).Value = 1
$ws.Cells.Item(192, 1).Value = 'Nepal'
$ws.Cells.Item(192, 4).Value = 1
$ws.Cells.Item(192, 8).Value = 0
$ws.Cells.Item(193, 1).Value = 'Gambia'
$ws.Cells.Item(194, 1).Value = 'Anguila'
$ws.Cells.Item(195, 1).Value = 'Belice'
$ws.Cells.Item(197, 1).Value = 'Islas Virgenes Britanicas'
$ws.Cells.Item(198, 1).Value = 'Somalia'
$ws.Cells.Item(199, 1).Value = 'Islas Turcas y Caicos'
$ws.Cells.Item(200, 1).Value = 'San Cristobal y Nieves'
$ws.Cells.Item(202, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(203, 1).Value = 'San Vicente y las Granadinas'
$ws.Cells.Item(204, 1).Value = 'Timor Oriental'
$ws.Cells.Item(205, 1).Value = 'Libia'
